$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shrink the data range from 7 rows to 5 rows and refresh the values ---
# Wipe the old A1:B7 block (drops rows 6 & 7 entirely) then write the new 5-row table.
$ws.Range("A1:B7").ClearContents()

$data = @(
    @(1, 265),
    @(2, 213),
    @(3, 224),
    @(4, 630),
    @(5, 615)
)
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# --- Update the bar chart to match: blank title, ranges trimmed to the 5 rows ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart

$chart.ChartTitle.Text = ""

$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(,Sheet1!`$A`$1:`$A`$5,Sheet1!`$B`$1:`$B`$5,1)"
